# Weekly update: insert two new price rows (most recent week) at the top of
# the data block (rows 204-205), pushing the existing historical rows down
# by two positions (old row 204 -> new row 206, ..., old row 252 -> new row 254).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 204; this shifts rows 204:252 down to 206:254
$ws.Rows("204:205").Insert()

# --- New row 204 ---
$ws.Range("A204").Value2 = 4
$ws.Range("B204").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C204").Value2 = "Los Lagos"
$ws.Range("D204").Value2 = 44641
$ws.Range("E204").Value2 = 10
$ws.Range("F204").Value2 = 100112024
$ws.Range("G204").Value2 = "Choclo"
$ws.Range("H204").Value2 = "Choclero"
$ws.Range("I204").Value2 = "Primera"
$ws.Range("J204").Value2 = 3000
$ws.Range("K204").Value2 = 250
$ws.Range("L204").Value2 = 300
$ws.Range("M204").Value2 = 275
$ws.Range("N204").Value2 = "`$/unidad"
$ws.Range("O204").Value2 = "Región del Maule"
$ws.Range("P204").Value2 = 275
$ws.Range("Q204").Value2 = 1
$ws.Range("R204").Value2 = "Hortaliza"

# --- New row 205 ---
$ws.Range("A205").Value2 = 4
$ws.Range("B205").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C205").Value2 = "Los Lagos"
$ws.Range("D205").Value2 = 44641
$ws.Range("E205").Value2 = 10
$ws.Range("F205").Value2 = 100112024
$ws.Range("G205").Value2 = "Choclo"
$ws.Range("H205").Value2 = "Dulce o Americano"
$ws.Range("I205").Value2 = "Primera"
$ws.Range("J205").Value2 = 6000
$ws.Range("K205").Value2 = 150
$ws.Range("L205").Value2 = 200
$ws.Range("M205").Value2 = 175
$ws.Range("N205").Value2 = "`$/unidad"
$ws.Range("O205").Value2 = "Región del Maule"
$ws.Range("P205").Value2 = 175
$ws.Range("Q205").Value2 = 1
$ws.Range("R205").Value2 = "Hortaliza"
